$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.028.67'
$ws.Range('E2').Value = '  -1.56%  '
$ws.Range('D3').Value = '3.521.62'
$ws.Range('E3').Value = '  -0.09%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = "'585.04"
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').Value = "'132.28"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.46%  '
$ws.Range('D7').Value = '3.522.25'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = "'0.485"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.32%  '
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('D11').Value = "'7.10"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.74%  '
$ws.Range('E12').Value = '  -2.47%  '
$ws.Range('D13').Value = '4.118.82'
$ws.Range('E13').Value = '  +0.08%  '
$ws.Range('D14').Value = "'27.60"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('E15').Value = '  +1.29%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.534.30'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value = "'0.0000179"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.81%  '
$ws.Range('D18').Value = '64.070.01'
$ws.Range('E18').Value = '  -1.38%  '
$ws.Range('D19').Value = "'9.85"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.10%  '
$ws.Range('E20').Value = '  -2.38%  '
$ws.Range('E21').Value = '  -1.70%  '
$ws.Range('D22').Value = "'384.67"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.02%  '
$ws.Range('E23').Value = '  -0.78%  '
$ws.Range('D24').Value = '3.660.35'
$ws.Range('E24').Value = '  -0.14%  '
$ws.Range('D25').Value = "'73.70"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.49%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').Value = "'0.0000116"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.57%  '
$ws.Range('D28').Value = "'1.58"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('D29').Value = "'7.51"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.69%  '
$ws.Range('E30').Value = '  -0.17%  '
$ws.Range('E31').Value = '  -0.79%  '
$ws.Range('E32').Value = '  -1.62%  '
$ws.Range('D33').Value = '3.530.24'
$ws.Range('E33').Value = '  +0.10%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').Value = "'23.54"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.47%  '
$ws.Range('E36').Value = '  +0.28%  '
$ws.Range('D37').Value = "'5.38"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.22%  '
$ws.Range('D38').Value = "'1.57"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.50%  '
$ws.Range('E39').Value = '  -1.18%  '
$ws.Range('D40').Value = "'158.79"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.80%  '
$ws.Range('D41').Value = "'0.0791"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.05%  '
$ws.Range('D42').Value = "'0.816"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.85%  '
$ws.Range('D43').Value = "'26.25"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.14%  '
$ws.Range('E44').Value = '  +0.28%  '
$ws.Range('D45').Value = "'41.86"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.51%  '
$ws.Range('E46').Value = '  -4.83%  '
$ws.Range('E47').Value = '  -0.65%  '
$ws.Range('E48').Value = '  -2.68%  '
$ws.Range('D49').Value = '2.444.13'
$ws.Range('E49').Value = '  +1.24%  '
$ws.Range('E50').Value = '  -1.47%  '
$ws.Range('D51').Value = "'0.910"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.26%  '
